# Insert a new weekly data row at row 122 (pushing the existing rows 122..225
# down to 123..226, growing the sheet's used range to A1:R226), then populate
# the newly inserted row with this week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(122).Insert()

$ws.Cells.Item(122, 1).Value = 7
$ws.Cells.Item(122, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(122, 3).Value = "Ñuble"
$ws.Cells.Item(122, 4).Value = 44589
$ws.Cells.Item(122, 5).Value = 16
$ws.Cells.Item(122, 6).Value = 100112023
$ws.Cells.Item(122, 7).Value = "Brócoli"
$ws.Cells.Item(122, 8).Value = "Sin especificar"
$ws.Cells.Item(122, 9).Value = "Primera"
$ws.Cells.Item(122, 10).Value = 300
$ws.Cells.Item(122, 11).Value = 700
$ws.Cells.Item(122, 12).Value = 750
$ws.Cells.Item(122, 13).Value = 725
$ws.Cells.Item(122, 14).Value = "$/unidad"
$ws.Cells.Item(122, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(122, 16).Value = 725
$ws.Cells.Item(122, 17).Value = 1
$ws.Cells.Item(122, 18).Value = "Hortaliza"

$ws.Cells.Item(122, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
